$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wall")

# Updated randomized batch data for rows 32-46
# Columns B:I = Run_ID, Wall, Roof, Glazing, WWR, Orient, Light, Equip
$rows = @(
    @(36, 0, 0, 0, 0, 0, 0, 1),
    @(37, 0, 0, 0, 0, 0, 0, 2),
    @(38, 0, 0, 0, 0, 0, 0, 3),
    @(39, 0, 0, 0, 0, 0, 0, 4),
    @(40, 0, 0, 0, 0, 0, 0, 5),
    @(46, 0, 0, 0, 0, 0, 1, 0),
    @(47, 0, 0, 0, 0, 0, 2, 0),
    @(48, 0, 0, 0, 0, 0, 3, 0),
    @(49, 0, 0, 0, 0, 0, 4, 0),
    @(50, 0, 0, 0, 0, 0, 5, 0),
    @(53, 0, 0, 1, 0, 0, 0, 0),
    @(54, 0, 0, 2, 0, 0, 0, 0),
    @(55, 0, 0, 3, 0, 0, 0, 0),
    @(56, 0, 0, 4, 0, 0, 0, 0),
    @(57, 0, 0, 5, 0, 0, 0, 0)
)

$arr = New-Object 'object[,]' $rows.Count,8
for ($i = 0; $i -lt $rows.Count; $i++) {
    for ($j = 0; $j -lt 8; $j++) {
        $arr[$i,$j] = $rows[$i][$j]
    }
}
$ws.Range("B32:I46").Value = $arr

# The old rows 37-39 and 47-49 blocks are no longer needed; remove the trailing rows
$ws.Rows("47:49").Delete()

# Restore the view selection (previously scrolled to A23 with S37 selected)
$ws.Range("N15").Select() | Out-Null
